$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LINE_TRIALS_URL")

# ------------------------------------------------------------------
# 1. Clear out the old trial rows (4-6) completely, and drop the
#    leftover E3 value (the old row 3 is being overwritten with new
#    trial data below, so only its stray E3 "TBD" cell needs clearing)
# ------------------------------------------------------------------
$ws.Range("A4:G6").Clear()
$ws.Range("E3").ClearContents()

# ------------------------------------------------------------------
# 2. Rename column header BOM_UNDER_TRIAL -> TRIAL
# ------------------------------------------------------------------
$ws.Range("B1").Value = "TRIAL"

# ------------------------------------------------------------------
# 3. New trial data (row 2 + row 3)
# ------------------------------------------------------------------
$tDate = (Get-Date -Year 2025 -Month 9 -Day 18).Date

$ws.Range("A2").Value = "JTPV"
$ws.Range("B2").Value = "R&D Production Order | G12R N-Type HEP Cell Line Trial | Comparison with`nregular 12R N-Type Cell"
$ws.Range("B2").WrapText = $true
$ws.Range("C2").Value = "In Progress"
$ws.Range("D2").Value = $tDate
$ws.Range("F2").Value = 23500

$ws.Range("A3").Value = "GCL, URECO and JTPV"
$ws.Range("B3").Value = "R&D Production Order | Line trial of N-type M10 cells"
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = "In Progress"
$ws.Range("D3").Value = $tDate
$ws.Range("F3").Value = 22420

# Row heights to match the wrapped-text content
$ws.Rows(2).RowHeight = 58
$ws.Rows(3).RowHeight = 29

# Widen the TRIAL column for the longer descriptions
$ws.Columns("B").ColumnWidth = 25.75

# ------------------------------------------------------------------
# 4. Move the lone formatted placeholder cell from K7 to J7 (it
#    shifts left by one once column G is dropped from the table)
# ------------------------------------------------------------------
$ws.Range("K7").Copy($ws.Range("J7"))
$ws.Range("K7").Clear()

# ------------------------------------------------------------------
# 5. Drop the trailing "Column1" table column and shrink the table /
#    autofilter range down from the full-sheet extent
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Item(7).Delete()
$lo.Resize($ws.Range("A1:F1048572"))

# ------------------------------------------------------------------
# 6. Rebuild data validations against the new (smaller) ranges
# ------------------------------------------------------------------
$ws.Range("A2:A1048576").Validation.Delete()
$ws.Range("B2:B1048576").Validation.Delete()
$ws.Range("C2:C1048576").Validation.Delete()
$ws.Range("D2:D1048576").Validation.Delete()
$ws.Range("E2").Validation.Delete()
$ws.Range("E20:E1048576").Validation.Delete()
$ws.Range("E3:E19").Validation.Delete()

$ws.Range("A2:A1048572").Validation.Add(3, 1, 1, "=Vendor_List")
$ws.Range("B2:B1048572").Validation.Add(3, 1, 1, "=Material_List")
$ws.Range("C2:C1048572").Validation.Add(3, 1, 1, "=Status_List")
$ws.Range("D2:D1048572").Validation.Add(4, 1, 1, "=36526", "=73050")
$ws.Range("E2").Validation.Add(4, 1, 1, "=36526", "=73050")
$ws.Range("E16:E1048576").Validation.Add(4, 1, 1, "=36526", "=73050")
$ws.Range("E3:E15").Validation.Add(0, 1, 1)

# ------------------------------------------------------------------
# 7. Update the active selection to reflect where the editor left off
# ------------------------------------------------------------------
$ws.Range("B6").Select()

Write-Host "done"
